$d = $word.ActiveDocument

$replacements = @(
    @("995÷2=", "482÷2="),
    @("354÷2=", "964÷4="),
    @("271÷2=", "685÷9="),
    @("956÷8=", "195÷8="),
    @("371÷4=", "773÷9="),
    @("671÷6=", "895÷2="),
    @("901÷6=", "789÷6="),
    @("496÷8=", "482÷2="),
    @("606÷6=", "869÷9="),
    @("752÷7=", "897÷6="),
    @("796÷2=", "207÷4="),
    @("753÷5=", "564÷7="),
    @("896÷7=", "158÷5="),
    @("691÷4=", "197÷2="),
    @("123÷5=", "377÷2="),
    @("699÷3=", "237÷3="),
    @("770÷6=", "560÷4="),
    @("615÷2=", "809÷8="),
    @("331÷8=", "882÷2="),
    @("408÷4=", "792÷6="),
    @("169÷8=", "904÷3="),
    @("373÷4=", "293÷3="),
    @("287÷5=", "542÷6="),
    @("148÷6=", "514÷6="),
    @("392÷3=", "453÷6="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                            $true, 1, $false, $new, 2) | Out-Null
}
